# Adjust cost files relying on US data
# - Add an "India Adjustment" worksheet (before BRAaCTSC) that computes a
#   multiplier from the ratio of India/US LDV gasoline prices.
# - Reference that multiplier from the BRAaCTSC sheet's cost formula.
# - Add an explanatory note on the About sheet.

$wb = $excel.ActiveWorkbook

# --- 1. About sheet: add explanatory note in A23 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A23").Value = "For India, we apply an adjustment factor based on the ratio of LDV prices in India and the US."

# --- 2. Insert new "India Adjustment" sheet right before BRAaCTSC ---
$wsBRAaCTSC = $wb.Worksheets.Item("BRAaCTSC")
$wsIndia = $wb.Worksheets.Add($wsBRAaCTSC)
$wsIndia.Name = "India Adjustment"

$wsIndia.Range("A1").Value = "India passenger LDV price, gasoline"
$wsIndia.Range("A2").Value = "US passenger LDV price, gasoline"
$wsIndia.Range("A3").Value = "Multiplier"

$wsIndia.Range("B1").Value = 7677.1720917226739
$wsIndia.Range("B2").Value = 36108.776943503013
$wsIndia.Range("B3").Formula = "=B1/B2"

$wsIndia.Range("D1").Value = "*see variable trans/BNVP, India EPS 3.1"
$wsIndia.Range("D2").Value = "*see variable trans/BNVP, US EPS 3.1"

# Target authored width is 33.59765625 characters; the host's ColumnWidth
# setter only resolves to 1/6-character increments, so 32.8 is the closest
# input that lands on the nearest achievable stored width (33.67).
$wsIndia.Columns.Item(1).ColumnWidth = 32.8

# --- 3. Update BRAaCTSC cost formula to apply the India multiplier ---
# Re-fetch the worksheet by name since sheet indices shifted after the insert.
$wsBRAaCTSC = $wb.Worksheets.Item("BRAaCTSC")
$wsBRAaCTSC.Range("B2").Formula = "=Data!B6*'India Adjustment'!B3"

# --- 4. Restore view state (selected cells / active sheet) ---
[void]$wsAbout.Range("A24").Select()
[void]$wsIndia.Range("D3").Select()
[void]$wsBRAaCTSC.Range("B3").Select()
$wsAbout.Activate()
